# regen s_vals data to filter save games
# Update columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-26.
# Column A (date) and F (Win) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 10.76772829162301)
    3  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    4  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    5  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    6  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    7  = @(1.455362044514542, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 12.1040582534506)
    8  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    9  = @(0.01293466051926884, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0.9634143985795411)
    10 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    11 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 13.26162679800665)
    12 = @(3.286832544864788, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 4.23731228292506)
    13 = @(0.2917716402565462, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.364486309189372)
    14 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    15 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    16 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    17 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548)
    18 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    19 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502)
    20 = @(0.2917716402565462, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0.9761466351224579)
    21 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    22 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    23 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    24 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
    25 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    26 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}

Write-Host "Updated s_vals rows 2-26 (B,C,D,E,G); filtered save games"
